$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("AI", "AJ", "AN")

for ($row = 2; $row -le 21; $row++) {
    foreach ($col in $columns) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        $val = $cell.Value()
        if ($val -eq $null -or $val -eq "") {
            $cell.Value = "/"
        }
    }
}
